$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update percent coverage values (deep sea double count fix)
$ws.Range("B2").Value = 66.67211941872959
$ws.Range("C2").Value = 88.76230889935026

$ws.Range("C3").Value = 92.31243281328341

$ws.Range("B4").Value = 84.16103557019441
$ws.Range("C4").Value = 91.36211995779584

$ws.Range("C9").Value = 95.33920971283109

$ws.Range("C10").Value = 95.35073559458699

$ws.Range("C11").Value = 98.48406223524952

$ws.Range("C12").Value = 75.5966977831971

$ws.Range("C14").Value = 98.26540129188768

$ws.Range("C15").Value = 87.55520630885273

$ws.Range("C16").Value = 88.25445330474211

# Update note text in D19 to mention 'Deep Sea' and remove ISSCAAP code 46
$noteText = "NOTE: Percent coverages are performed across FAO major fishing areas to be consistent with Fishstatj. `nThus, landings from areas such as 'Salmon', 'Tuna', 'Deep Sea', and 'Sharks' are added back into the FAO major fishing area from where they were reported. `nPercent coverage calculations do not include landings from ISSCAAP codes 61, 62, 63, 64, 71, 72, 73, 74, 81, 82, 83, 91, 92, 93, 94, `nexcept for stocks from these groups which are included in the assessment."
$ws.Range("D19").Value = $noteText
